# Add show "Your Friendly Neighborhood Spider-Man" and "Iron Man and His Awesome Friends"
# to the "animated-tv-series" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("animated-tv-series")

# Row 35: Your Friendly Neighborhood Spider-Man
$ws.Range("A35").Value = "animated-tv-series"
$ws.Range("B35").Value = "2025"
$ws.Range("C35").Value = "Your Friendly Neighborhood Spider-Man"
$ws.Range("D35").Value = "https://en.wikipedia.org/wiki/Your_Friendly_Neighborhood_Spider-Man"
$ws.Range("F35").Value = 1
$ws.Range("G35").Value = 10

# Row 36: Iron Man and His Awesome Friends
$ws.Range("A36").Value = "animated-tv-series"
$ws.Range("B36").Value = "2025"
$ws.Range("C36").Value = "Iron Man and His Awesome Friends"
$ws.Range("D36").Value = "https://en.wikipedia.org/wiki/Iron_Man_and_His_Awesome_Friends"
$ws.Range("F36").Value = 1
$ws.Range("G36").Value = 10

# Match the style (text format) used on the rest of column B.
$ws.Range("B35").NumberFormat = "@"
$ws.Range("B36").NumberFormat = "@"

$ws.Range("H37").Select()
